$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.777.84"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.926.38"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'603.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "'166.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("D7").Value = "3.921.77"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").Value = "'6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'37.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "4.579.24"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "3.904.28"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").Value = "68.871.66"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "'7.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").Value = "'17.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "'11.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "'484.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").Value = "'0.0000170"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.94%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "'0.721"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "'84.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "'12.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").Value = "4.073.43"
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "'32.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "3.872.11"
$ws.Range("E35").Value = "  +3.33%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "'5.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "'3.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("D42").Value = "'0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").Value = "'433.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'8.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'27.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.31%  "
$ws.Range("D49").Value = "2.836.46"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'141.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +17.24%  "
